$d = $word.ActiveDocument

# Locate the paragraph that contains the "Ver no Jupiter..." text and the
# paragraph right before it (an empty "Normal" paragraph) as well as the
# paragraph right after it ("(c) 2020 ... Creative Commons Attribution").
# All three of these paragraphs (the blank spacer, the "Ver no Jupiter..."
# line, and the copyright line) are removed, while the blank paragraph
# that originally followed the copyright line, and everything else, is
# left untouched.

$count = $d.Paragraphs.Count
$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $startIndex = $i - 1
    }
    if ($text -like "*Creative Commons Attribution*") {
        $endIndex = $i
    }
}

if ($startIndex -ge 1 -and $endIndex -ge $startIndex) {
    $rangeStart = $d.Paragraphs.Item($startIndex).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIndex).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
